$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143342322"
$ws.Range("D16").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E16").Value = "2104"
$ws.Range("F16").Value = 69334
$ws.Range("G16").Value = 2000000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143342322"
$ws.Range("D17").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E17").Value = "2103"
$ws.Range("F17").Value = 80000
$ws.Range("G17").Value = 2000000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143342322"
$ws.Range("D18").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E18").Value = "2102"
$ws.Range("F18").Value = 80000
$ws.Range("G18").Value = 2000000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143342322"
$ws.Range("D19").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E19").Value = "2101"
$ws.Range("F19").Value = 80000
$ws.Range("G19").Value = 2000000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143342322"
$ws.Range("D20").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E20").Value = "2012"
$ws.Range("F20").Value = 80000
$ws.Range("G20").Value = 2000000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143342322"
$ws.Range("D21").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E21").Value = "2011"
$ws.Range("F21").Value = 80000
$ws.Range("G21").Value = 2000000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1143342322"
$ws.Range("D22").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E22").Value = "2010"
$ws.Range("F22").Value = 80000
$ws.Range("G22").Value = 2000000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1143342322"
$ws.Range("D23").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E23").Value = "2009"
$ws.Range("F23").Value = 80000
$ws.Range("G23").Value = 2000000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1143342322"
$ws.Range("D24").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E24").Value = "2008"
$ws.Range("F24").Value = 80000
$ws.Range("G24").Value = 2000000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1143342322"
$ws.Range("D25").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E25").Value = "2007"
$ws.Range("F25").Value = 80000
$ws.Range("G25").Value = 2000000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1143342322"
$ws.Range("D26").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E26").Value = "2006"
$ws.Range("F26").Value = 80000
$ws.Range("G26").Value = 2000000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1143342322"
$ws.Range("D27").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E27").Value = "2005"
$ws.Range("F27").Value = 80000
$ws.Range("G27").Value = 2000000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1143342322"
$ws.Range("D28").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E28").Value = "2004"
$ws.Range("F28").Value = 80000
$ws.Range("G28").Value = 2000000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1143342322"
$ws.Range("D29").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E29").Value = "2003"
$ws.Range("F29").Value = 80000
$ws.Range("G29").Value = 2000000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "1143342322"
$ws.Range("D30").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E30").Value = "2002"
$ws.Range("F30").Value = 80000
$ws.Range("G30").Value = 2000000

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "1143342322"
$ws.Range("D31").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E31").Value = "2001"
$ws.Range("F31").Value = 80000
$ws.Range("G31").Value = 2000000

$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "1143342322"
$ws.Range("D32").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E32").Value = "1912"
$ws.Range("F32").Value = 80000
$ws.Range("G32").Value = 2000000

$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "1143342322"
$ws.Range("D33").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E33").Value = "1911"
$ws.Range("F33").Value = 80000
$ws.Range("G33").Value = 2000000

$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "1143342322"
$ws.Range("D34").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E34").Value = "1910"
$ws.Range("F34").Value = 80000
$ws.Range("G34").Value = 2000000

$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "1143342322"
$ws.Range("D35").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E35").Value = "1909"
$ws.Range("F35").Value = 80000
$ws.Range("G35").Value = 2000000

$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "1143342322"
$ws.Range("D36").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E36").Value = "1908"
$ws.Range("F36").Value = 80000
$ws.Range("G36").Value = 2000000

$ws.Range("B37").Value = "CC"
$ws.Range("C37").Value = "1143342322"
$ws.Range("D37").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E37").Value = "1907"
$ws.Range("F37").Value = 80000
$ws.Range("G37").Value = 2000000

$ws.Range("B38").Value = "CC"
$ws.Range("C38").Value = "1143342322"
$ws.Range("D38").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E38").Value = "1906"
$ws.Range("F38").Value = 80000
$ws.Range("G38").Value = 2000000

$ws.Range("B39").Value = "CC"
$ws.Range("C39").Value = "1143342322"
$ws.Range("D39").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E39").Value = "1905"
$ws.Range("F39").Value = 80000
$ws.Range("G39").Value = 2000000

$ws.Range("B40").Value = "CC"
$ws.Range("C40").Value = "1143342322"
$ws.Range("D40").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E40").Value = "1904"
$ws.Range("F40").Value = 80000
$ws.Range("G40").Value = 2000000

$ws.Range("B41").Value = "CC"
$ws.Range("C41").Value = "1143342322"
$ws.Range("D41").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E41").Value = "1903"
$ws.Range("F41").Value = 80000
$ws.Range("G41").Value = 2000000

$ws.Range("B42").Value = "CC"
$ws.Range("C42").Value = "1143342322"
$ws.Range("D42").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E42").Value = "1902"
$ws.Range("F42").Value = 80000
$ws.Range("G42").Value = 2000000

$ws.Range("B43").Value = "CC"
$ws.Range("C43").Value = "1143342322"
$ws.Range("D43").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E43").Value = "1901"
$ws.Range("F43").Value = 80000
$ws.Range("G43").Value = 2000000

$ws.Range("B44").Value = "CC"
$ws.Range("C44").Value = "1143342322"
$ws.Range("D44").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E44").Value = "1812"
$ws.Range("F44").Value = 80000
$ws.Range("G44").Value = 2000000

$ws.Range("B45").Value = "CC"
$ws.Range("C45").Value = "1143342322"
$ws.Range("D45").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E45").Value = "1811"
$ws.Range("F45").Value = 80000
$ws.Range("G45").Value = 2000000

$ws.Range("B46").Value = "CC"
$ws.Range("C46").Value = "1143342322"
$ws.Range("D46").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E46").Value = "1810"
$ws.Range("F46").Value = 80000
$ws.Range("G46").Value = 2000000

$ws.Range("B47").Value = "CC"
$ws.Range("C47").Value = "1143342322"
$ws.Range("D47").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E47").Value = "1809"
$ws.Range("F47").Value = 80000
$ws.Range("G47").Value = 2000000

$ws.Range("B48").Value = "CC"
$ws.Range("C48").Value = "1143342322"
$ws.Range("D48").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E48").Value = "1808"
$ws.Range("F48").Value = 80000
$ws.Range("G48").Value = 2000000

$ws.Range("B49").Value = "CC"
$ws.Range("C49").Value = "1143342322"
$ws.Range("D49").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E49").Value = "1807"
$ws.Range("F49").Value = 80000
$ws.Range("G49").Value = 2000000

$ws.Range("B50").Value = "CC"
$ws.Range("C50").Value = "1143342322"
$ws.Range("D50").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E50").Value = "1806"
$ws.Range("F50").Value = 80000
$ws.Range("G50").Value = 2000000

$ws.Range("B51").Value = "CC"
$ws.Range("C51").Value = "1143342322"
$ws.Range("D51").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E51").Value = "1805"
$ws.Range("F51").Value = 80000
$ws.Range("G51").Value = 2000000

$ws.Range("B52").Value = "CC"
$ws.Range("C52").Value = "1143342322"
$ws.Range("D52").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E52").Value = "1804"
$ws.Range("F52").Value = 80000
$ws.Range("G52").Value = 2000000

$ws.Range("B53").Value = "CC"
$ws.Range("C53").Value = "1143342322"
$ws.Range("D53").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E53").Value = "1803"
$ws.Range("F53").Value = 80000
$ws.Range("G53").Value = 2000000

$ws.Range("B54").Value = "CC"
$ws.Range("C54").Value = "1143342322"
$ws.Range("D54").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E54").Value = "1802"
$ws.Range("F54").Value = 80000
$ws.Range("G54").Value = 2000000

$ws.Range("B55").Value = "CC"
$ws.Range("C55").Value = "1143342322"
$ws.Range("D55").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E55").Value = "1801"
$ws.Range("F55").Value = 80000
$ws.Range("G55").Value = 2000000

$ws.Range("B56").Value = "CC"
$ws.Range("C56").Value = "1143342322"
$ws.Range("D56").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E56").Value = "1712"
$ws.Range("F56").Value = 80000
$ws.Range("G56").Value = 2000000

$ws.Range("B57").Value = "CC"
$ws.Range("C57").Value = "1143342322"
$ws.Range("D57").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E57").Value = "1711"
$ws.Range("F57").Value = 80000
$ws.Range("G57").Value = 2000000

$ws.Range("B58").Value = "CC"
$ws.Range("C58").Value = "1143342322"
$ws.Range("D58").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E58").Value = "1710"
$ws.Range("F58").Value = 80000
$ws.Range("G58").Value = 2000000

$ws.Range("B59").Value = "CC"
$ws.Range("C59").Value = "1143342322"
$ws.Range("D59").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E59").Value = "1709"
$ws.Range("F59").Value = 80000
$ws.Range("G59").Value = 2000000

$ws.Range("B60").Value = "CC"
$ws.Range("C60").Value = "1143342322"
$ws.Range("D60").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E60").Value = "1708"
$ws.Range("F60").Value = 80000
$ws.Range("G60").Value = 2000000

$ws.Range("B61").Value = "CC"
$ws.Range("C61").Value = "1143342322"
$ws.Range("D61").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E61").Value = "1707"
$ws.Range("F61").Value = 80000
$ws.Range("G61").Value = 2000000

$ws.Range("B62").Value = "CC"
$ws.Range("C62").Value = "1143342322"
$ws.Range("D62").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E62").Value = "1706"
$ws.Range("F62").Value = 80000
$ws.Range("G62").Value = 2000000

$ws.Range("B63").Value = "CC"
$ws.Range("C63").Value = "1143342322"
$ws.Range("D63").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E63").Value = "1705"
$ws.Range("F63").Value = 80000
$ws.Range("G63").Value = 2000000

$ws.Range("B64").Value = "CC"
$ws.Range("C64").Value = "1143342322"
$ws.Range("D64").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E64").Value = "1703"
$ws.Range("F64").Value = 80000
$ws.Range("G64").Value = 2000000

$ws.Range("B65").Value = "CC"
$ws.Range("C65").Value = "1143342322"
$ws.Range("D65").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E65").Value = "1702"
$ws.Range("F65").Value = 80000
$ws.Range("G65").Value = 2000000

$ws.Range("B66").Value = "CC"
$ws.Range("C66").Value = "1143342322"
$ws.Range("D66").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E66").Value = "1701"
$ws.Range("F66").Value = 80000
$ws.Range("G66").Value = 2000000

$ws.Range("B67").Value = "CC"
$ws.Range("C67").Value = "1143342322"
$ws.Range("D67").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E67").Value = "1612"
$ws.Range("F67").Value = 80000
$ws.Range("G67").Value = 2000000

$ws.Range("B68").Value = "CC"
$ws.Range("C68").Value = "1143342322"
$ws.Range("D68").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E68").Value = "1611"
$ws.Range("F68").Value = 80000
$ws.Range("G68").Value = 2000000

$ws.Range("B69").Value = "CC"
$ws.Range("C69").Value = "1143342322"
$ws.Range("D69").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E69").Value = "1610"
$ws.Range("F69").Value = 80000
$ws.Range("G69").Value = 2000000

$ws.Range("B70").Value = "CC"
$ws.Range("C70").Value = "1143342322"
$ws.Range("D70").Value = "SHIRLEY YOHANA ROMERO VELANDIA"
$ws.Range("E70").Value = "1609"
$ws.Range("F70").Value = 80000
$ws.Range("G70").Value = 2000000

$ws.Range("B71").Value = "CC"
$ws.Range("C71").Value = "72313674"
$ws.Range("D71").Value = "VELINTON JASITH FONTALVO ARIZA"
$ws.Range("E71").Value = "1712"
$ws.Range("F71").Value = 80000
$ws.Range("G71").Value = 2000000

$ws.Range("B72").Value = "CC"
$ws.Range("C72").Value = "72313674"
$ws.Range("D72").Value = "VELINTON JASITH FONTALVO ARIZA"
$ws.Range("E72").Value = "1711"
$ws.Range("F72").Value = 80000
$ws.Range("G72").Value = 2000000

$ws.Range("B73").Value = "CC"
$ws.Range("C73").Value = "72313674"
$ws.Range("D73").Value = "VELINTON JASITH FONTALVO ARIZA"
$ws.Range("E73").Value = "1710"
$ws.Range("F73").Value = 80000
$ws.Range("G73").Value = 2000000

$ws.Range("B74").Value = "CC"
$ws.Range("C74").Value = "72313674"
$ws.Range("D74").Value = "VELINTON JASITH FONTALVO ARIZA"
$ws.Range("E74").Value = "1709"
$ws.Range("F74").Value = 80000
$ws.Range("G74").Value = 2000000

$ws.Range("B75").Value = "CC"
$ws.Range("C75").Value = "72313674"
$ws.Range("D75").Value = "VELINTON JASITH FONTALVO ARIZA"
$ws.Range("E75").Value = "1708"
$ws.Range("F75").Value = 80000
$ws.Range("G75").Value = 2000000
